$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata": bump version/date, fill in Publisher, replace the
# duplicated "Contact" row with a single "Jurisdiction" row.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail"; row 11 duplicated it.
# Delete row 11 (the duplicate) so the table shrinks from 21 to 20 rows, then
# turn what is now the single remaining row into "Jurisdiction".
$meta.Rows.Item(11).Delete()

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---------------------------------------------------------------------------
# Sheet "Elements": update the Short/Definition text of the root Extension
# row from the generic placeholder to the parameter-value-specific text.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Parameter Value"
$elements.Range("L2").Value = "Value that was used for the parameter"
